$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestLoginLogout")

# Rename sheet1
$ws.Name = "OrangeTC1"

# Set column widths (values chosen so the resulting stored width is as close
# as possible to the target 15.5703125 / 16.140625 / 14.5703125 given the
# engine's internal pixel-quantization of ColumnWidth)
$ws.Columns.Item(1).ColumnWidth = 14.666666666666666
$ws.Columns.Item(2).ColumnWidth = 15.333333333333334
$ws.Columns.Item(3).ColumnWidth = 13.666666666666666

# Update existing cells
$ws.Range("A2").Value = "opensourcecms"
$ws.Range("B2").Value = "opensourcecms"

# Add new header cells
$ws.Range("C1").Value = "Homepageurl"
$ws.Range("D1").Value = "Loginpageurl"

# Add new data cells (order matters for shared string table indices)
$ws.Range("D2").Value = "login"
$ws.Range("C2").Value = "index.php"

# Set selection to C2
$ws.Range("C2").Select()
